# Insert a new data row at row 6 (pushing existing rows 6..71 down to 7..72)
# and populate it with the new weekly price record, per the commit
# "Fruta / hortaliza, semanal" (weekly fruit/vegetable price update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 6:71 down by inserting a new row at position 6.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new record.
$ws.Cells.Item(6, 1).Value  = 9
$ws.Cells.Item(6, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(6, 3).Value  = "Metropolitana"
$ws.Cells.Item(6, 4).Value  = 44496
$ws.Cells.Item(6, 5).Value  = 13
$ws.Cells.Item(6, 6).Value  = 100112022
$ws.Cells.Item(6, 7).Value  = "Arveja Verde"
$ws.Cells.Item(6, 8).Value  = "Perfection"
$ws.Cells.Item(6, 9).Value  = "Primera"
$ws.Cells.Item(6, 10).Value = 30
$ws.Cells.Item(6, 11).Value = 24000
$ws.Cells.Item(6, 12).Value = 24000
$ws.Cells.Item(6, 13).Value = 24000
$ws.Cells.Item(6, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(6, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(6, 16).Value = 960
$ws.Cells.Item(6, 17).Value = 25
$ws.Cells.Item(6, 18).Value = "Hortaliza"
